# "added new endpoints and refactor"
#
# A new meter reading (dated 7/6/2024) is the latest entry and is inserted
# as the new row 2, pushing the previous row 2 (dated 7/5/2024) down to
# row 3 and renumbering it from 1 -> 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data row (row 2) down to make room for the new reading.
$ws.Rows.Item(2).Insert()

# New row 2: latest reading.
# Value2 with a leading apostrophe keeps numeric/date-looking text
# ("1", "12345678", "7/6/2024") stored as text instead of being coerced
# to a number or date serial by Excel's smart-entry parsing.
$ws.Cells.Item(2, 1).Value2 = "'1"
$ws.Cells.Item(2, 2).Value2 = "'d2f89a58-4961-4fa3-baad-3f7eff79ce02"
$ws.Cells.Item(2, 3).Value2 = "'12345678"
$ws.Cells.Item(2, 4).Value2 = "'7/6/2024"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 2

# Row 3 (the original row 2, shifted down by the insert): renumber it.
$ws.Cells.Item(3, 1).Value2 = "'2"

# The insert carried the header row's bold/centered style onto row 2;
# restore both data rows to the workbook's plain default style.
$ws.Range("A2:G3").Style = "Normal"
